# Add 2022-Q4 data
# ------------------------------------------------------------------
# 1) Insert a brand-new worksheet named "2022-Q4" right after "总计"
#    (i.e. right before the current first quarterly sheet "2022-Q3").
# 2) Populate that new sheet with the fund-holding table for 2022-Q4.
# 3) Update the "总计" (summary) sheet: push the existing quarter rows
#    down by one row and insert the new 2022-Q4 summary row at the top.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

# ---- Step 1: create the new "2022-Q4" sheet right after "总计" ----
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# ---- Step 2: populate the new sheet, matching the layout/style of ----
# ---- the other quarterly sheets (copy header/format from 2022-Q3) ----
$usedCols = 8
for ($c = 1; $c -le $usedCols; $c++) {
    $q3Sheet.Cells.Item(1, $c).Copy()
    $q4Sheet.Cells.Item(1, $c).PasteSpecial(-4122)
}
for ($c = 1; $c -le $usedCols; $c++) {
    $q3Sheet.Cells.Item(2, $c).Copy()
    $q4Sheet.Cells.Item(2, $c).PasteSpecial(-4122)
}

$q4Sheet.Cells.Item(1, 2).Value = "基金代码"
$q4Sheet.Cells.Item(1, 3).Value = "基金名称"
$q4Sheet.Cells.Item(1, 4).Value = "基金规模"
$q4Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q4Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q4Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4Sheet.Cells.Item(1, 8).Value = "仓位排名"

$q4Rows = @(
    @(0, "011868", "中信建投远见回报混合A", "6.34", "95.01", "3.44", "0.2181", 9),
    @(1, "011869", "中信建投远见回报混合C", "1.41", "95.01", "3.44", "0.0485", 9),
    @(2, "710002", "富安达策略精选混合",     "0.58", "67.44", "2.05", "0.0119", 5)
)

$r = 2
foreach ($row in $q4Rows) {
    if ($r -gt 2) {
        for ($c = 1; $c -le $usedCols; $c++) {
            $q3Sheet.Cells.Item($r, $c).Copy()
            $q4Sheet.Cells.Item($r, $c).PasteSpecial(-4122)
        }
    }
    $q4Sheet.Cells.Item($r, 1).Value = $row[0]
    $q4Sheet.Cells.Item($r, 2).Value = $row[1]
    $q4Sheet.Cells.Item($r, 3).Value = $row[2]
    $q4Sheet.Cells.Item($r, 4).Value = $row[3]
    $q4Sheet.Cells.Item($r, 5).Value = $row[4]
    $q4Sheet.Cells.Item($r, 6).Value = $row[5]
    $q4Sheet.Cells.Item($r, 7).Value = $row[6]
    $q4Sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---- Step 3: update the "总计" sheet ----
# Push rows 2-6 down to rows 3-7 (copy formats + values from bottom up
# so we never overwrite data we still need to read). Column A is just a
# zero-based running index (0,1,2,...) so it is re-derived, not copied.
for ($row = 6; $row -ge 2; $row--) {
    $dst = $row + 1
    for ($c = 1; $c -le 4; $c++) {
        $totalSheet.Cells.Item($row, $c).Copy()
        $totalSheet.Cells.Item($dst, $c).PasteSpecial(-4122)
    }
    $totalSheet.Cells.Item($dst, 2).Value = $totalSheet.Cells.Item($row, 2).Value()
    $totalSheet.Cells.Item($dst, 3).Value = $totalSheet.Cells.Item($row, 3).Value()
    $totalSheet.Cells.Item($dst, 4).Value = $totalSheet.Cells.Item($row, 4).Value()
    $totalSheet.Cells.Item($dst, 1).Value = $dst - 2
}

# Now write the brand-new 2022-Q4 summary row into row 2.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.28

Write-Host "2022-Q4 sheet added and 总计 sheet updated"
